$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.511.48'
$ws.Range("E2").Value = '  -1.13%  '
$ws.Range("D3").Value = '1.851.66'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6526'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07509'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2983'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.50'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07635'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("D12").Value = '1.852.44'
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("E13").Value = '  -0.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6852'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.80'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009497'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.121'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '29.535.99'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").Value = '2.121.62'
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.24'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9995'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.691'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.34%  '
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("E27").Value = '  -0.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.81'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06037'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.488'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.45%  '
$ws.Range("E31").Value = '  -0.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.143'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.071'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.73%  '
$ws.Range("E34").Value = '  +1.55%  '
$ws.Range("E35").Value = '  -0.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7243'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.589'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.803'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01782'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.75%  '
$ws.Range("D40").Value = '1.201.99'
$ws.Range("E40").Value = '  -2.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.245'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9094'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9994'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("E44").Value = '  -1.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.89'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.50'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("E47").Value = '  +10.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000123'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.92%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4058'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.125'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.660'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.12%  '
